# Fruta / hortaliza, semanal
# Swap the data of rows 2 & 4, and rows 3 & 5 (columns D, K, L, M, N, O, P, Q, S)
# to reflect the updated weekly price report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell($row1, $row2, $col) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value = $v2
    $ws.Range($addr2).Value = $v1
}

$columns = @('D', 'K', 'L', 'M', 'N', 'O', 'P', 'Q', 'S')

foreach ($col in $columns) {
    Swap-Cell 2 4 $col
    Swap-Cell 3 5 $col
}
